$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Tarantula")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 3.14355481300603
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 10.55712633337616

$ws = $wb.Worksheets.Item("Ochiai")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.113610075825722
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 9.996144454440273

$ws = $wb.Worksheets.Item("Op2")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 4.424881120678577
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 13.45392623056161

$ws = $wb.Worksheets.Item("Barinel")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 3.146125176712494
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 10.55969669708263

$ws = $wb.Worksheets.Item("Dstar")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.096902711733704
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 10.83665338645415

$ws = $wb.Worksheets.Item("Russell_rao")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 7.161033286210009
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 19.86634108726387

$ws = $wb.Worksheets.Item("Simple_matching")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 2.770852075568659
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 16.57820331576913

$ws = $wb.Worksheets.Item("Rogers_tanimoto")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 2.770852075568659
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 16.57820331576913

$ws = $wb.Worksheets.Item("Ample")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.445186993959635
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 13.81699010409972

$ws = $wb.Worksheets.Item("Jaccard")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.122606348798347
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 9.739750674720442

$ws = $wb.Worksheets.Item("Cohen")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.121963757871731
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 6.765197275414458

$ws = $wb.Worksheets.Item("Scott")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.121321166945115
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 6.754273229661985

$ws = $wb.Worksheets.Item("Rogot1")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.121321166945115
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 6.754273229661985

$ws = $wb.Worksheets.Item("Geometric_mean")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.113610075825722
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 7.038941010152916

$ws = $wb.Worksheets.Item("M2")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 2.175170286595544
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 15.74797583858113

$ws = $wb.Worksheets.Item("Wong1")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 7.161033286210009
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 19.86634108726387

$ws = $wb.Worksheets.Item("Sokal")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 2.770852075568659
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 16.57820331576913

$ws = $wb.Worksheets.Item("Sorensen_dice")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.122606348798347
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 9.739750674720442

$ws = $wb.Worksheets.Item("Dice")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.122606348798347
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 9.739750674720442

$ws = $wb.Worksheets.Item("Humman")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 2.770852075568659
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 16.57820331576913

$ws = $wb.Worksheets.Item("Wong2")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 2.770852075568659
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 16.57820331576913

$ws = $wb.Worksheets.Item("Euclid")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 2.770852075568659
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 16.57820331576913

$ws = $wb.Worksheets.Item("Zoltar")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 2.831255622670597
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 12.92443130702994

$ws = $wb.Worksheets.Item("Rogot2")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.143169258450063
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 10.68435933684615

$ws = $wb.Worksheets.Item("Hamming")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 2.770852075568659
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 16.57820331576913

$ws = $wb.Worksheets.Item("Fleiss")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.109754530266025
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 7.900012851818524

$ws = $wb.Worksheets.Item("Anderberg")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.122606348798347
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 9.739750674720442

$ws = $wb.Worksheets.Item("Goodman")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.122606348798347
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 9.739750674720442

$ws = $wb.Worksheets.Item("Harmonic_mean")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.143169258450063
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 7.977123763012458

$ws = $wb.Worksheets.Item("Kulczynski2")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.319239172342874
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 12.68410230047551
